$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (139:140), pushing the
# existing rows 139-226 down to 141-228 (and growing the used range to T228).
$ws.Rows("139:140").Insert()

# Row 139: new entry for "1a nueva(o)" quality, Provincia de Quillota origin.
$ws.Range("A139").Value = 4
$ws.Range("B139").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C139").Value = "Los Lagos"
$ws.Range("D139").Value = 44488
$ws.Range("E139").Value = 10
$ws.Range("F139").Value = "Fruta"
$ws.Range("G139").Value = 100106
$ws.Range("H139").Value = "Oleaginosos"
$ws.Range("I139").Value = 100106002
$ws.Range("J139").Value = "Palta"
$ws.Range("K139").Value = "Hass"
$ws.Range("L139").Value = "1a nueva(o)"
$ws.Range("M139").Value = 300
$ws.Range("N139").Value = 4000
$ws.Range("O139").Value = 4200
$ws.Range("P139").Value = 4100
$ws.Range("Q139").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R139").Value = "Provincia de Quillota"
$ws.Range("S139").Value = 4100
$ws.Range("T139").Value = 1

# Row 140: new entry for "2a nueva(o)" quality, Provincia de Quillota origin.
$ws.Range("A140").Value = 4
$ws.Range("B140").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C140").Value = "Los Lagos"
$ws.Range("D140").Value = 44488
$ws.Range("E140").Value = 10
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100106
$ws.Range("H140").Value = "Oleaginosos"
$ws.Range("I140").Value = 100106002
$ws.Range("J140").Value = "Palta"
$ws.Range("K140").Value = "Hass"
$ws.Range("L140").Value = "2a nueva(o)"
$ws.Range("M140").Value = 150
$ws.Range("N140").Value = 3600
$ws.Range("O140").Value = 3600
$ws.Range("P140").Value = 3600
$ws.Range("Q140").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R140").Value = "Provincia de Quillota"
$ws.Range("S140").Value = 3600
$ws.Range("T140").Value = 1
